$d = $word.ActiveDocument

# --- Update the 5 existing paragraphs in place ---

# 1. Update date in the daily header line (17.06.24 -> 16.06.24)
$ok1 = $d.Content.Find.Execute('⚡️🚀המאמר היומי של מייק 17.06.24:⚡️🚀', $true, $false, $false, $false, $false, $true, 1, $false, '⚡️🚀המאמר היומי של מייק 16.06.24:⚡️🚀', 2)

# 2. Update paper title (SSAMBA -> STATISTICAL REJECTION SAMPLING ...)
$ok2 = $d.Content.Find.Execute('SSAMBA: SELF-SUPERVISED AUDIO REPRESENTATION LEARNING WITH MAMBA STATE SPACE MODEL', $true, $false, $false, $false, $false, $true, 1, $false, 'STATISTICAL REJECTION SAMPLING IMPROVES PREFERENCE OPTIMIZATION', 2)

# 3. Replace short intro line with the new long intro paragraph
# (uses Find/Replace so the internal triple-space run text is not
#  marked xml:space="preserve" by the writer, matching the target)
$ok3 = $d.Content.Find.Execute('הסקירה נמצאת כאן: ', $true, $false, $false, $false, $false, $true, 1, $false, 'המאמר הזה וכמה הבאים שאסקור בימים הקרובים מציעים שכלולים שונים לשיטה Direct Preference Optimization או בקיצור DPO. למעשה DPO בעצמה היא שדרוג של Proximal Policy Optimization או PPO שהפכה להיות מאוד פופולרית אחרי שמכמה חברות השתמשו בה ליישור מודלי שפה (alignment   או instruction tuning) בתור השלב האחרון של אימון מודל שפה foundational. השיטה שייכת למשפחת RLHF כי היא דורשת דאטה (שאלות ותשובות) המדורגות על ידי בני אדם - עבור כל שאלה הם (המתייגים) בוחרים מה התשובה איזה תשובה טובה יותר.', 2)

# 4. Replace google docs review link with DPO explanation paragraph
# (direct Range.Text assignment keeps the literal straight quote " intact;
#  Find/Replace's replacement text would get smart-quoted by AutoCorrect)
$d.Paragraphs.Item(4).Range.Text = 'למעשה DPO בא לייתר את מודל התגמול (reward) גם חוסך גם משאבים לאימונו וגם מאפשר לא להחזיק מודל נוסף בשלב RLHF. למעשה DPO מנצל את המבנה של פונקצית לוס של PPO, שהיא מקסום פונקציית תגמול עם איבר רגולריזציה שבא לשמור את המודל המיושר קרוב למודל התחלתי, כדי להיפטר מפונקציית התגמול בפונקציית לוס. זה מתאפשר עקב העובדה שקיים ביטוי מפורש לפוליסי האופטימלי (מודל שפה ״מושלם אחרי היישור") דרך הפוליסי אחרי ה-SFT (מודל שפה שאנו מתחילים ממנו את אימון היישור) ופונקציית התגמול.'

# 5. Replace old arxiv link (2405.11831) with new explanatory paragraph text
$ok5 = $d.Content.Find.Execute('https://arxiv.org/abs/2405.11831', $true, $false, $false, $false, $false, $true, 1, $false, 'אחרי שמשתמשים במודל לוס המושרה על ידי מודל (Bradley-Terry (BT המגדיר מהי הסתברות העדפה של תשובה חיובית על תשובה שלילית (על אותה השאלה) מה- rewards שלהם, ואנו מגיעים לביטוי עבור לוס של RLHF שמכיל רק את הפוליסי התחלתי. זה למעשה DPO והוא ממזער את פונקציית הלוס שלו על סט המכיל  זוגות של תשובות טובות וגרועות.', 2)

# --- Append 8 brand-new paragraphs after the (former link, now 5th) paragraph ---
$baseIndex = 5
$d.Paragraphs.Item($baseIndex + 0).Range.InsertParagraphAfter()
$d.Paragraphs.Item($baseIndex + 1).Range.Text = 'המאמר שנסקור היום שואל את השאלה האם הדגימה האחידה מהסט הזה היא אופטימלית (מבחינת איכות התוצאה שהיא הופליסי הסופי או מודל שפה אחרי היישור). אולי אם היה לנו פונקציית תגמול היינו מעדיפים זוגות עם יחס מקסימלי בין ה-reward של התשובה החיובית לשלולית? אולי צריך לתעדף זוגות עם reward שלילי הנמוך ביותר?'
$d.Paragraphs.Item($baseIndex + 1).Range.InsertParagraphAfter()
$d.Paragraphs.Item($baseIndex + 2).Range.Text = 'המאמר מציע את הגישה הבאה:'
$d.Paragraphs.Item($baseIndex + 2).Range.InsertParagraphAfter()
$d.Paragraphs.Item($baseIndex + 3).Range.Text = '- מאמנים מודל text2text שבהינתן שאלה ושתי תשובות מוציא את התשובה המועדפת.'
$d.Paragraphs.Item($baseIndex + 3).Range.InsertParagraphAfter()
$d.Paragraphs.Item($baseIndex + 4).Range.Text = '- בעזרת המודל הזה בונים את פונקציית התגמול דרך סמלוץ (על ידי דגימה של שאלה וזוג תשובות) של הסתברות העדפה של תשובה טובה על תשובה גרועה.'
$d.Paragraphs.Item($baseIndex + 4).Range.InsertParagraphAfter()
$d.Paragraphs.Item($baseIndex + 5).Range.Text = '- בעזרת פונקציית תגמול זו בונים פוליסי pi_r שלמעשה זה מודל שפה (המאפשר לחשב הסתברות של תשובה בהינתן שאלה)'
$d.Paragraphs.Item($baseIndex + 5).Range.InsertParagraphAfter()
$d.Paragraphs.Item($baseIndex + 6).Range.Text = '- משתמשים בדגימת rejection כדי לדגום pi_r באמצעות הפוליסי ההתחלתי (= מודל שפה) כדי למזער את הלוס בדרך לפוליסי "המיושר".'
$d.Paragraphs.Item($baseIndex + 6).Range.InsertParagraphAfter()
$d.Paragraphs.Item($baseIndex + 7).Range.Text = 'הם גם משחקים עם כמה פונקציות לוס כמו hinge loss (בטח כבר שכחתם אבל אוהבים להשתמש בו ב -SVM).'
$d.Paragraphs.Item($baseIndex + 7).Range.InsertParagraphAfter()
$d.Paragraphs.Item($baseIndex + 8).Range.Text = 'https://arxiv.org/abs/2309.06657'

Write-Output ("FindResults: " + $ok1 + "," + $ok2 + "," + $ok3 + "," + $ok5)
Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
